$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, pushing existing rows 24..54 down to 25..55.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new weekly price record.
$ws.Cells.Item(24,1).Value = 10
$ws.Cells.Item(24,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(24,3).Value = "La Araucanía"
$ws.Cells.Item(24,4).Value = "11/28/2022"
$ws.Cells.Item(24,5).Value = 9
$ws.Cells.Item(24,6).Value = 100112042
$ws.Cells.Item(24,7).Value = "Locoto"
$ws.Cells.Item(24,8).Value = "Sin especificar"
$ws.Cells.Item(24,9).Value = "Primera"
$ws.Cells.Item(24,10).Value = 50
$ws.Cells.Item(24,11).Value = 2500
$ws.Cells.Item(24,12).Value = 2500
$ws.Cells.Item(24,13).Value = 2500
$ws.Cells.Item(24,14).Value = "`$/kilo"
$ws.Cells.Item(24,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(24,16).Value = 2500
$ws.Cells.Item(24,17).Value = 1
$ws.Cells.Item(24,18).Value = "Hortaliza"
